$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-28
# from serial date 45221 (2023-10-22) to 45224 (2023-10-25)
for ($row = 2; $row -le 28; $row++) {
    $ws.Cells.Item($row, 3).Value = 45224
}
